$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "2024-06-02" (4th sheet): fix up rows 11-13, drop row 14.
#    Each patient now gets their own distinct doctor visit record, so the
#    duplicated rows are replaced with the corrected data and the extra
#    trailing row is removed.
# ---------------------------------------------------------------------------
$ws0602 = $wb.Worksheets.Item("2024-06-02")

$ws0602.Range("B11").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws0602.Range("B11").Value = 45445.86062358796
$ws0602.Range("E11").Value = "2006-05-29"
$ws0602.Range("F11").Value = "porrkthso[rhk"
$ws0602.Range("G11").Value = "435678"

$ws0602.Range("B12").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws0602.Range("B12").Value = 45445.860871875
$ws0602.Range("E12").Value = "2006-05-29"
$ws0602.Range("F12").Value = "dkjghsropjk"
$ws0602.Range("G12").Value = "456879"

$ws0602.Range("B13").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws0602.Range("B13").Value = 45445.87829535879
$ws0602.Range("C13").Value = "jhkjhblijho"
$ws0602.Range("E13").Value = "2006-05-29"
$ws0602.Range("F13").Value = "fyckuj"

$ws0602.Rows.Item(14).Delete()

# ---------------------------------------------------------------------------
# 2) Two new daily-log sheets get appended: "2024-07-09" and "2024-07-10".
#    Their layout now carries "Врач" / "Врач_Индекс" columns per visit.
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws0709 = $wb.Worksheets.Add($null, $afterSheet)
$ws0709.Name = "2024-07-09"

$ws0709.Range("A1").Value = "ID"
$ws0709.Range("B1").Value = "Время"
$ws0709.Range("C1").Value = "ФИО пациента"
$ws0709.Range("D1").Value = "М\Ж\Р"
$ws0709.Range("E1").Value = "Дата рождения"
$ws0709.Range("F1").Value = "Причина"
$ws0709.Range("G1").Value = "Давление"
$ws0709.Range("H1").Value = "Врач"
$ws0709.Range("I1").Value = "Врач_Индекс"

$ws0709.Range("A2").NumberFormat = "@"
$ws0709.Range("A2").Value = "1"
$ws0709.Range("B2").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws0709.Range("B2").Value = 45482.79670672453
$ws0709.Range("C2").Value = "шортшотш"
$ws0709.Range("D2").Value = "Ж"
$ws0709.Range("E2").Value = "2006-04-06"
$ws0709.Range("F2").Value = "рототщто"
$ws0709.Range("G2").Value = "7890"

$ws0709.Columns.Item(2).ColumnWidth = 17.81640625

$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws0710 = $wb.Worksheets.Add($null, $afterSheet2)
$ws0710.Name = "2024-07-10"

$ws0710.Range("A1").Value = "ID"
$ws0710.Range("B1").Value = "Время"
$ws0710.Range("C1").Value = "ФИО пациента"
$ws0710.Range("D1").Value = "Врач"
$ws0710.Range("E1").Value = "Врач_Индекс"
$ws0710.Range("F1").Value = "М\Ж\Р"
$ws0710.Range("G1").Value = "Дата рождения"
$ws0710.Range("H1").Value = "Причина"
$ws0710.Range("I1").Value = "Давление"

$ws0710.Range("A2").NumberFormat = "@"
$ws0710.Range("A2").Value = "1"
$ws0710.Range("B2").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws0710.Range("B2").Value = 45483.5933899537
$ws0710.Range("C2").Value = "ааонег"
$ws0710.Range("D2").Value = "Karp_Kuzmin"
$ws0710.Range("E2").Value = 3
$ws0710.Range("F2").Value = "Ж"
$ws0710.Range("G2").Value = "2006-07-04"
$ws0710.Range("H2").Value = "шкгпщш"
$ws0710.Range("I2").Value = "7890"

$ws0710.Range("A3").NumberFormat = "@"
$ws0710.Range("A3").Value = "2"
$ws0710.Range("B3").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws0710.Range("B3").Value = 45483.59384497685
$ws0710.Range("C3").Value = "квгнгпрщщж"
$ws0710.Range("D3").Value = "Karp_Kuzmin"
$ws0710.Range("E3").Value = 3
$ws0710.Range("F3").Value = "М"
$ws0710.Range("G3").Value = "2006-07-03"
$ws0710.Range("H3").Value = "dstfui"
$ws0710.Range("I3").Value = "ytfu67589"

$ws0710.Range("A4").NumberFormat = "@"
$ws0710.Range("A4").Value = "3"
$ws0710.Range("B4").NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws0710.Range("B4").Value = 45483.59436456019
$ws0710.Range("C4").Value = "рпплродж"
$ws0710.Range("D4").Value = "Karp_Kuzmin"
$ws0710.Range("E4").Value = 3
$ws0710.Range("F4").Value = "М"
$ws0710.Range("G4").Value = "2006-07-03"
$ws0710.Range("H4").Value = "utfuygu"
$ws0710.Range("I4").Value = "678"

$ws0710.Columns.Item(2).ColumnWidth = 17.81640625
$ws0710.Columns.Item(4).ColumnWidth = 11.7265625

# ---------------------------------------------------------------------------
# 3) The "current" summary sheet drops the doctor-name / doctor-id columns
#    (every visit now has its own doctor, so a single "doctor of the day"
#    column no longer makes sense) and gains the two new days.
#    Columns B ("ФИО врача") and, after the shift, C ("id_doctor") are
#    removed outright so the surviving columns keep their original widths
#    (bestFit etc.) instead of being recreated from scratch.
# ---------------------------------------------------------------------------
$wsCurrent = $wb.Worksheets.Item("current")
$wsCurrent.Columns.Item(2).Delete()
$wsCurrent.Columns.Item(3).Delete()

# Only the two totals that actually changed need correcting; the rest
# survived the column shift with their original values intact.
$wsCurrent.Range("B3").Value = 12
$wsCurrent.Range("D3").Value = 6

# Append the two new days.
$wsCurrent.Range("A4:A5").NumberFormat = "@"

$wsCurrent.Range("A4").Value = "2024-07-09"
$wsCurrent.Range("B4").Value = 1
$wsCurrent.Range("C4").Value = 0
$wsCurrent.Range("D4").Value = 1
$wsCurrent.Range("E4").Value = 0

$wsCurrent.Range("A5").Value = "2024-07-10"
$wsCurrent.Range("B5").Value = 3
$wsCurrent.Range("C5").Value = 0
$wsCurrent.Range("D5").Value = 1
$wsCurrent.Range("E5").Value = 2

$wsCurrent.Range("B7").Select()

# ---------------------------------------------------------------------------
# 4) Make the newest day sheet the active tab, matching the "current work"
#    focus moving to the latest log.
# ---------------------------------------------------------------------------
$ws0710.Activate()
